$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# Add a new row of data to the bottom of the table (expands Table1 ref/autoFilter
# from A1:E26 to A1:E27, and inherits formatting from the row above)
$newRow = $lo.ListRows.Add()

$ws.Range("A27").Value = 45599
$ws.Range("B27").Value = "Spez 1"
$ws.Range("C27").Value = "Zaubertrank"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 3

# Select new cell where user would end up after adding a row (A28)
$ws.Range("A28").Select()
